$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '60.307.77'
$ws.Cells.Item(3, 4).Value = '3.298.52'
$ws.Cells.Item(3, 5).Value = '  -3.68%  '
$ws.Cells.Item(4, 5).Value = '  +0.04%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '557.71'
$ws.Cells.Item(5, 5).Value = '  -3.87%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '140.62'
$ws.Cells.Item(6, 5).Value = '  -8.91%  '
$ws.Cells.Item(7, 5).Value = '  -0.03%  '
$ws.Cells.Item(8, 4).Value = '3.299.62'
$ws.Cells.Item(8, 5).Value = '  -3.63%  '
$ws.Cells.Item(9, 5).Value = '  -3.65%  '
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '7.91'
$ws.Cells.Item(10, 5).Value = '  -2.10%  '
$ws.Cells.Item(11, 5).Value = '  -5.33%  '
$ws.Cells.Item(12, 5).Value = '  -2.83%  '
$ws.Cells.Item(13, 4).Value = '3.863.32'
$ws.Cells.Item(13, 5).Value = '  -3.73%  '
$ws.Cells.Item(14, 5).Value = '  -0.36%  '
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '26.53'
$ws.Cells.Item(15, 5).Value = '  -7.69%  '
$ws.Cells.Item(16, 4).Value = '3.298.34'
$ws.Cells.Item(16, 5).Value = '  -3.66%  '
$ws.Cells.Item(17, 5).Value = '  -5.10%  '
$ws.Cells.Item(18, 4).Value = '60.289.86'
$ws.Cells.Item(18, 5).Value = '  -2.98%  '
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '6.07'
$ws.Cells.Item(19, 5).Value = '  -7.31%  '
$ws.Cells.Item(20, 5).Value = '  -5.58%  '
$ws.Cells.Item(21, 5).Value = '  -5.49%  '
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '373.42'
$ws.Cells.Item(22, 5).Value = '  -2.74%  '
$ws.Cells.Item(23, 5).Value = '  -0.02%  '
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '72.10'
$ws.Cells.Item(24, 5).Value = '  -5.04%  '
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '0.531'
$ws.Cells.Item(25, 5).Value = '  -7.19%  '
$ws.Cells.Item(26, 4).Value = '3.432.44'
$ws.Cells.Item(26, 5).Value = '  -3.66%  '
$ws.Cells.Item(27, 5).Value = '  -10.12%  '
$ws.Cells.Item(28, 5).Value = '  -1.71%  '
$ws.Cells.Item(29, 5).Value = '  +0.33%  '
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '7.02'
$ws.Cells.Item(30, 5).Value = '  -8.57%  '
$ws.Cells.Item(31, 5).Value = '  -0.05%  '
$ws.Cells.Item(32, 5).Value = '  -5.13%  '
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '7.30'
$ws.Cells.Item(33, 5).Value = '  -7.46%  '
$ws.Cells.Item(34, 5).Value = '  -3.28%  '
$ws.Cells.Item(35, 5).Value = '  -5.11%  '
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '165.53'
$ws.Cells.Item(36, 5).Value = '  -1.81%  '
$ws.Cells.Item(37, 5).Value = '  -9.10%  '
$ws.Cells.Item(38, 5).Value = '  -5.22%  '
$ws.Cells.Item(39, 5).Value = '  -5.30%  '
$ws.Cells.Item(40, 4).Value = '3.328.46'
$ws.Cells.Item(40, 5).Value = '  -3.83%  '
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '0.0722'
$ws.Cells.Item(41, 5).Value = '  -7.74%  '
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '25.40'
$ws.Cells.Item(42, 5).Value = '  -17.77%  '
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '41.81'
$ws.Cells.Item(43, 5).Value = '  -2.31%  '
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '0.746'
$ws.Cells.Item(44, 5).Value = '  -4.46%  '
$ws.Cells.Item(45, 5).Value = '  -4.37%  '
$ws.Cells.Item(46, 5).Value = '  -7.48%  '
$ws.Cells.Item(47, 5).Value = '  -6.67%  '
$ws.Cells.Item(48, 5).Value = '  +0.01%  '
$ws.Cells.Item(49, 4).Value = '2.320.83'
$ws.Cells.Item(49, 5).Value = '  -9.28%  '
$ws.Cells.Item(50, 5).Value = '  -6.53%  '
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '21.43'
$ws.Cells.Item(51, 5).Value = '  -8.68%  '
